{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same changes as the target diff:\n//  - Bold the 5 non-\"\u65e5\u671f\" header cells of the trends table.\n//  - Reword 3 header cells (\"...\u603b\u9500\u91cf\" -> \"...\u603b\u9500\u552e\u6570\u91cf\", etc.) and the\n//    \"\u5728\u7ebf\u641c\u7d22\u91cf\" header cell.\n//  - Normalize several date cells in the first (date) column.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// --- Header row (row 0): bold + text updates ---------------------------\nconst headerUpdates = [\n  { col: 1, text: \"\u5370\u5ea6\u5976\u8336\u603b\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n  { col: 2, text: \"\u624b\u5de5\u5370\u5ea6\u5976\u8336\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n  { col: 3, text: \"\u9884\u5236\u5370\u5ea6\u5976\u8336\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n  { col: 4, text: null },\n  { col: 5, text: \"\u5370\u5ea6\u5976\u8336\u7684\u5728\u7ebf\u641c\u7d22\u60c5\u51b5\" }\n];\n\nfor (const { col, text } of headerUpdates) {\n  const cell = table.getCell(0, col);\n  if (text !== null) cell.value = text;\n  cell.body.font.bold = true;\n}\n\n// --- Date column (col 0): normalize specific row values -----------------\nconst dateUpdates = [\n  { row: 3, text: \"2023/2/28\" },\n  { row: 6, text: \"2023/5/31\" },\n  { row: 8, text: \"2023/7/30\" },\n  { row: 9, text: \"2023/8/29\" },\n  { row: 10, text: \"2023/9/30\" },\n  { row: 12, text: \"2023/11/30\" }\n];\n\nfor (const { row, text } of dateUpdates) {\n  table.getCell(row, 0).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same changes as the target diff:\n#  - Bold the 5 non-\"\u65e5\u671f\" header cells of the trends table.\n#  - Reword 3 header cells (\"...\u603b\u9500\u91cf\" -> \"...\u603b\u9500\u552e\u6570\u91cf\", etc.) and the\n#    \"\u5728\u7ebf\u641c\u7d22\u91cf\" header cell.\n#  - Normalize several date cells in the first (date) column.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Header row (row 1): bold + text updates ----------------------------\n$headerUpdates = @(\n    @{ Col = 2; Text = \"\u5370\u5ea6\u5976\u8336\u603b\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n    @{ Col = 3; Text = \"\u624b\u5de5\u5370\u5ea6\u5976\u8336\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n    @{ Col = 4; Text = \"\u9884\u5236\u5370\u5ea6\u5976\u8336\u9500\u552e\u6570\u91cf\uff08\u4ef6\uff09\" },\n    @{ Col = 5; Text = $null },\n    @{ Col = 6; Text = \"\u5370\u5ea6\u5976\u8336\u7684\u5728\u7ebf\u641c\u7d22\u60c5\u51b5\" }\n)\n\nforeach ($u in $headerUpdates) {\n    $cell = $t.Cell(1, $u.Col)\n    if ($u.Text -ne $null) {\n        $cell.Range.Text = $u.Text\n    }\n    $cell.Range.Font.Bold = 1\n}\n\n# --- Date column (col 1): normalize specific row values ------------------\n$dateUpdates = @(\n    @{ Row = 4;  Text = \"2023/2/28\" },\n    @{ Row = 7;  Text = \"2023/5/31\" },\n    @{ Row = 9;  Text = \"2023/7/30\" },\n    @{ Row = 10; Text = \"2023/8/29\" },\n    @{ Row = 11; Text = \"2023/9/30\" },\n    @{ Row = 13; Text = \"2023/11/30\" }\n)\n\nforeach ($u in $dateUpdates) {\n    $t.Cell($u.Row, 1).Range.Text = $u.Text\n}\n"}
